# Update the workbook per the recorded diff:
#  1. Column C ("Förändrad" / last-changed date) moves from serial 45184
#     to serial 45186 for every data row (rows 2-262).
#  2. The HYPERLINK() formulas in columns S, T, U, V, W, X, Y (only present
#     on rows 2-10) gain a second HYPERLINK argument: the row's "Beteckning"
#     (column A) value, used as the friendly display text for the link.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 262

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 3).Value = 45186
}

$linkCols = @("S", "T", "U", "V", "W", "X", "Y")

for ($row = $firstRow; $row -le 10; $row++) {
    $beteckning = $ws.Cells.Item($row, 1).Value2

    foreach ($col in $linkCols) {
        $cell = $ws.Range($col + $row)
        $formula = $cell.Formula
        if ([string]::IsNullOrEmpty($formula)) { continue }
        if ($formula.TrimEnd().EndsWith(")")) {
            $newFormula = $formula.Substring(0, $formula.Length - 1) + ', "' + $beteckning + '")'
            $cell.Formula = $newFormula
        }
    }
}
